$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: the "Price" column (D) holds text-formatted numbers (e.g. "546.15").
# Values that parse as plain numbers must be forced to Text format first,
# otherwise Excel auto-converts them to numeric cells on assignment.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "60.156.84"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.320.20"
$ws.Range("E3").Value = "  -1.27%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.15"
$ws.Range("E5").Value = "  +0.08%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.15"
$ws.Range("E6").Value = "  -1.42%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -2.03%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.317.79"
$ws.Range("E9").Value = "  -1.25%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.07%  "

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.60"
$ws.Range("E11").Value = "  +1.58%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.45%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  +0.30%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.56"
$ws.Range("E14").Value = "  -1.25%  "

# Row 15 / Row 16 - WrappedBTC and WrappedliquidstakedEther2.0 swap places
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.734.07"
$ws.Range("E15").Value = "  -1.18%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "60.185.47"
$ws.Range("E16").Value = "  +0.16%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +0.64%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.323.06"
$ws.Range("E18").Value = "  -1.78%  "

# Row 19 - Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.55"
$ws.Range("E19").Value = "  -1.31%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  -1.67%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.40"
$ws.Range("E21").Value = "  -0.18%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -3.57%  "

# Row 23 - Dai
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.10%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.28"
$ws.Range("E24").Value = "  +1.64%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  -0.26%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.19%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.83"
$ws.Range("E27").Value = "  -0.74%  "

# Row 28 - Fetch.AI
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.36"
$ws.Range("E28").Value = "  +0.40%  "

# Row 29 - SuiNetwork
$ws.Range("E29").Value = "  +9.50%  "

# Row 30 - Monero
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.18"
$ws.Range("E30").Value = "  -0.31%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.19%  "

# Row 32 - PEPE
$ws.Range("E32").Value = "  +0.09%  "

# Row 33 - Aptos
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.04"
$ws.Range("E33").Value = "  +1.77%  "

# Row 34 - PolygonEcosystemToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.382"
$ws.Range("E34").Value = "  -0.14%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  -3.83%  "

# Row 36 - EthereumClassic
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.97"
$ws.Range("E36").Value = "  -0.44%  "

# Row 37 - USDe
$ws.Range("E37").Value = "  +0.00%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  -0.10%  "

# Row 39 - NEARProtocol
$ws.Range("E39").Value = "  -2.08%  "

# Row 40 - Bittensor
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "317.02"
$ws.Range("E40").Value = "  -0.68%  "

# Row 41 - OKB
$ws.Range("E41").Value = "  -0.30%  "

# Row 42 - Stacks
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.51"

# Row 43 - Aave
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.68"
$ws.Range("E43").Value = "  -2.92%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  +1.20%  "

# Row 45 - Stellar
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0944"
$ws.Range("E45").Value = "  -0.58%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.12"
$ws.Range("E46").Value = "  -1.46%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  +0.61%  "

# Row 48 - Hedera
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0494"
$ws.Range("E48").Value = "  -0.76%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +0.42%  "

# Row 50 - BabyDogeCoin
$ws.Range("D50").Value = "0.0₆0218"
$ws.Range("E50").Value = "  +2.12%  "

# Row 51 - WhiteBITCoin
$ws.Range("E51").Value = "  -0.88%  "
